$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.320.94'
$ws.Range("E2").Value = '  -3.49%  '

$ws.Range("D3").Value = '3.517.60'
$ws.Range("E3").Value = '  -5.03%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.80'
$ws.Range("E5").Value = '  -1.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.89'
$ws.Range("E6").Value = '  -3.80%  '

$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").Value = '3.510.01'
$ws.Range("E8").Value = '  -5.05%  '

$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("E10").Value = '  -6.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.76'
$ws.Range("E11").Value = '  +4.16%  '

$ws.Range("E12").Value = '  -2.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.08'
$ws.Range("E13").Value = '  -6.00%  '

$ws.Range("E14").Value = '  -3.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '676.69'
$ws.Range("E15").Value = '  -1.09%  '

$ws.Range("D16").Value = '4.083.20'
$ws.Range("E16").Value = '  -5.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.75'
$ws.Range("E17").Value = '  -3.32%  '

$ws.Range("D18").Value = '69.352.49'
$ws.Range("E18").Value = '  -3.46%  '

$ws.Range("D19").Value = '3.521.71'
$ws.Range("E19").Value = '  -4.76%  '

$ws.Range("E20").Value = '  -1.34%  '

$ws.Range("E21").Value = '  -3.73%  '

$ws.Range("E22").Value = '  -4.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.905'
$ws.Range("E23").Value = '  -4.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.24'
$ws.Range("E24").Value = '  -9.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.14'
$ws.Range("E25").Value = '  -5.65%  '

$ws.Range("E26").Value = '  -4.19%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.67'
$ws.Range("E28").Value = '  -6.17%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.47'
$ws.Range("E29").Value = '  -8.31%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.05'
$ws.Range("E30").Value = '  -7.20%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.77'
$ws.Range("E31").Value = '  -5.75%  '

$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.20'
$ws.Range("E32").Value = '  -7.84%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.36'
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.36'
$ws.Range("E34").Value = '  -5.90%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '598.45'
$ws.Range("E35").Value = '  +6.35%  '

$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.62'
$ws.Range("E36").Value = '  -15.43%  '

$ws.Range("B37").Value = 'Cosmos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.90'
$ws.Range("E37").Value = '  -3.66%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("E38").Value = '  -5.07%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.32'
$ws.Range("E39").Value = '  -4.01%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.27%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0440'
$ws.Range("E41").Value = '  -5.36%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.337'
$ws.Range("E42").Value = '  -4.52%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.424.88'
$ws.Range("E43").Value = '  -8.92%  '

$ws.Range("E44").Value = '  -6.23%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.46'
$ws.Range("E45").Value = '  -6.49%  '

$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0710'
$ws.Range("E46").Value = '  -9.35%  '

$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  +0.61%  '

$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.61'
$ws.Range("E48").Value = '  -7.20%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.80'
$ws.Range("E50").Value = '  +18.14%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.11'
$ws.Range("E51").Value = '  -2.14%  '

Write-Host "done"